# Update workbook for data through 2022-08-19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-08-11"

# Update the August row label (header text in column A)
$ws.Range("A9").Value = "August (through 08-11)"

# Update August row (row 9) values
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 21
$ws.Range("D9").Value = 24
$ws.Range("E9").Value = 18
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 74
$ws.Range("H9").Value = 72
$ws.Range("I9").Value = 61

# Update Total row (row 10) values
$ws.Range("B10").Value = 173
$ws.Range("C10").Value = 323
$ws.Range("D10").Value = 489
$ws.Range("E10").Value = 443
$ws.Range("F10").Value = 319
$ws.Range("G10").Value = 695
$ws.Range("H10").Value = 982
$ws.Range("I10").Value = 1031
